# The workbook's category "ScreenRecStarted" was renamed to "0_unstated".
# That label only ever appears as literal text (no formulas) in the header
# cell G1, and as part of the concatenated row-label strings in A27:A30
# (rowCategory + colCategory). Updating these five cells causes Excel to
# drop the now-unused "ScreenRecStarted" shared string and rebuild the
# shared-strings table accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "0_unstated"
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Move/record the active selection, matching the saved sheet view.
$ws.Range("E14").Select()
